$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# --- Column H width now matches the other narrow columns ---
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# --- Crime-data numbers (Week/28-Day/YTD/2-Yr/historical % changes) ---
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -90.909090909090
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 14
$ws.Range("I16").Value = 19
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = 5.555555555555
$ws.Range("L16").Value = 26.666666666666
$ws.Range("M16").Value = -17.391304347826
$ws.Range("N16").Value = -91.203703703703
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 23.809523809523
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 26
$ws.Range("K17").Value = 23.076923076923
$ws.Range("L17").Value = -8.571428571428
$ws.Range("M17").Value = 3.225806451612
$ws.Range("N17").Value = -54.285714285714
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -64.705882352941
$ws.Range("M18").Value = -75
$ws.Range("N18").Value = -98.165137614678
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -9.375
$ws.Range("I19").Value = 33
$ws.Range("J19").Value = 36
$ws.Range("K19").Value = -8.333333333333
$ws.Range("L19").Value = -46.774193548387
$ws.Range("M19").Value = -40
$ws.Range("N19").Value = -58.75
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -58.333333333333
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = -68.421052631578
$ws.Range("N20").Value = -97.5
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -8.045977011494
$ws.Range("I21").Value = 97
$ws.Range("J21").Value = 101
$ws.Range("K21").Value = -3.960396039603
$ws.Range("L21").Value = -31.690140845070
$ws.Range("M21").Value = -37.012987012987
$ws.Range("N21").Value = -89.778714436248
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -17.241379310344
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 146
$ws.Range("H24").Value = -33.561643835616
$ws.Range("I24").Value = 122
$ws.Range("J24").Value = 165
$ws.Range("K24").Value = -26.060606060606
$ws.Range("L24").Value = -18.666666666666
$ws.Range("M24").Value = 8.928571428571
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 91
$ws.Range("H25").Value = -29.670329670329
$ws.Range("I25").Value = 82
$ws.Range("J25").Value = 101
$ws.Range("K25").Value = -18.811881188118
$ws.Range("L25").Value = 6.493506493506
$ws.Range("C26").Value = 8
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 44
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = -12
$ws.Range("L26").Value = -2.222222222222
$ws.Range("M26").Value = -37.142857142857
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = -75
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = 0
$ws.Range("M29").Value = -100
$ws.Range("M30").Value = -100

# --- Cells that became "no activity" placeholders ("0" / "***.*") ---
# Source cells already holding the shared placeholder strings with the right style:
#   C14 -> "0"     E14 -> "***.*"
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("F33"))
